$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 59 (shifts the existing rows 59:76 down to 60:77)
$ws.Rows("59:59").Insert()

# Populate the newly inserted row with the new weekly observation
$ws.Range("A59").Value = 8
$ws.Range("B59").Value = "Terminal La Palmera de La Serena"
$ws.Range("C59").Value = "Coquimbo"
$ws.Range("D59").Value = 44875
$ws.Range("E59").Value = 4
$ws.Range("F59").Value = 100114007
$ws.Range("G59").Value = "Jengibre"
$ws.Range("H59").Value = "Sin especificar"
$ws.Range("I59").Value = "Primera"
$ws.Range("J59").Value = 480
$ws.Range("K59").Value = 14000
$ws.Range("L59").Value = 15000
$ws.Range("M59").Value = 14500
$ws.Range("N59").Value = "$/caja 13 kilos"
$ws.Range("O59").Value = "Perú"
$ws.Range("P59").Value = 1115
$ws.Range("Q59").Value = 13
$ws.Range("R59").Value = "Hortaliza"
